$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.173.32"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.685.11"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'215.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'23.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.82%  "
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("D10").Value = "'0.0626"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.928.07"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.693.61"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "'0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "'67.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "27.181.31"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'235.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'8.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "0.0₃0742"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("D23").Value = "'9.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.22%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "'147.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "'16.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "1.541.58"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "'3.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "'0.605"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").Value = "'0.944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").Value = "'2.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").Value = "'0.0173"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'69.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "1.834.24"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").Value = "'90.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.82%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").Value = "'8.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.24%  "
$ws.Range("E51").Value = "  +0.13%  "
